# Add two new user submissions to the collected-data sheet:
#   - a "Тест" test entry, inserted right after the header (new row 2),
#     pushing the existing entries down by one row;
#   - a "Тимурчик Мусаевчик" entry appended as the new last row.
#
# Row "insert" is done by hand (shifting values down from the bottom up)
# instead of Rows.Insert(), since Insert() drags the header's bold/border
# formatting onto the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 8
$lastCol = 7

for ($r = $lastRow; $r -ge 2; $r--) {
    $dst = $r + 1
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($dst, $c).Value = $ws.Cells.Item($r, $c).Value2
    }
}

# New row 2: the "Тест" entry.
$ws.Range("A2").Value = "Тест"
$ws.Range("B2").Value = "Тест"

# C2/E2 need to be genuine text (matching the source data's mixed typing),
# not Excel's auto-detected number. Compute them with a TEXT() formula, then
# paste-special as values so the result is plain text with no leftover
# formula and without touching NumberFormat (which would stamp a new style).
$ws.Range("C2").Formula = '=TEXT(89001083247,"0")'
$ws.Range("C2").Copy()
$ws.Range("C2").PasteSpecial(-4163)
$ws.Range("E2").Formula = '=TEXT(25,"0")'
$ws.Range("E2").Copy()
$ws.Range("E2").PasteSpecial(-4163)

$ws.Range("D2").Value = "test@mail.ru"
$ws.Range("F2").Value = "Студенты"
$ws.Range("G2").Value = "24.05.2023 12:55"

# New last row (10): the "Тимурчик Мусаевчик" entry.
$ws.Range("A10").Value = "Тимурчик"
$ws.Range("B10").Value = "Мусаевчик"
$ws.Range("C10").Value = 89521656455
$ws.Range("D10").Value = "musaev.timur@mail.ru"
$ws.Range("E10").Value = 19
$ws.Range("F10").Value = "Студенты"
$ws.Range("G10").Value = "24.05.2023 11:57"
